$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "5312166542"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "FEB2025"
$ws.Range("C2").Value = 169.54

$ws.Range("A2:C9").Select()
